$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format A4 as text first so the date-like string "2025-08-19" is stored
# literally instead of being auto-converted into a date serial number,
# then reset the style back to Normal so no extra formatting sticks.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-08-19"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 59.59000015258789
$ws.Range("C4").Value = 700.25
$ws.Range("D4").Value = 321.4500122070312
